$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that get "Yes" across columns B:K (default), with a few "No" overrides
$rowsWithData = @(2, 3, 4, 5, 6, 9, 10, 20, 22, 29)
$rowsWithNoInH = @(4, 5, 29)

foreach ($r in $rowsWithData) {
    $ws.Range("B$r`:K$r").Value = "Yes"
}

foreach ($r in $rowsWithNoInH) {
    $ws.Range("H$r").Value = "No"
}

$ws.Range("J12").Select()
